$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Force the new rows to be stored as text, matching the existing inline-string data in the sheet
$ws.Range("A17:K22").NumberFormat = "@"

$newRows = @(
    @("16", "11", "3", "0,0001", "784", "32", "10", "(0, 0,1)", "2,2077215917241975", "2,1794148056623137", "23,685714285714287"),
    @("17", "13", "3", "0,0001", "784", "32", "10", "(0, 0,1)", "2,223611367075916", "2,2048816641098052", "23,5"),
    @("18", "15", "3", "0,0001", "784", "32", "10", "(0, 0,1)", "2,2146582486515953", "2,1954575391601954", "23,714285714285715"),
    @("19", "17", "3", "0,0001", "784", "32", "10", "(0, 0,1)", "2,235946882498052", "2,223201451371017", "17,685714285714287"),
    @("20", "19", "3", "0,0001", "784", "32", "10", "(0, 0,1)", "2,2511713964947746", "2,2385759327792862", "17,214285714285715"),
    @("21", "1", "5", "0,0001", "784", "32", "10", "(0, 0,1)", "0,40360077396405736", "0,37413761358399433", "89,62857142857143")
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
